$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: only the Taxonsorteringsordning (B) value changes
$ws.Range("B4").Value = 90800

# Rows 5 and 6 swap their species/location/time data (columns A, E, F, G, H, Q, R, Z, AB),
# while column B (Taxonsorteringsordning) gets new independent values.

# --- Row 5 becomes the former row 6 data ---
$ws.Range("A5").Value = 112127587
$ws.Range("B5").Value = 90794
$ws.Range("E5").Value = 4362
$ws.Range("F5").Value = "Blå taggsvamp"
$ws.Range("G5").Value = "Hydnellum caeruleum"
$ws.Range("H5").Value = "(Hornem.) P.Karst."
$ws.Range("Q5").Value = 690447
$ws.Range("R5").Value = 7125629
$ws.Range("Z5").Value = "14:27"
$ws.Range("AB5").Value = "14:27"

# --- Row 6 becomes the former row 5 data ---
$ws.Range("A6").Value = 112127546
$ws.Range("B6").Value = 90792
$ws.Range("E6").Value = 4361
$ws.Range("F6").Value = "Orange taggsvamp"
$ws.Range("G6").Value = "Hydnellum aurantiacum"
$ws.Range("H6").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q6").Value = 690408
$ws.Range("R6").Value = 7125570
$ws.Range("Z6").Value = "14:25"
$ws.Range("AB6").Value = "14:25"
